# Apply the edits described by the diff:
# - The shared text "968-MS-EI-DB-DL-REC-NON-RNI-CTRFD-DL-MD-TR-1-Late Repayment"
#   becomes "968-MS-EI-DB-DL-REC-NON-RNI-CTRFD-DL-MD-TR-1-LateRepayment" (space removed)
#   and is used as the value of B1 on both sheets, now styled like the rest of column B.
# - The active sheet switches from "ProductLoanInput" to "ProductLoanOutput".
# - Selections reset to B1 on both sheets.

$wb = $excel.ActiveWorkbook

$wsIn = $wb.Worksheets.Item("ProductLoanInput")
$wsOut = $wb.Worksheets.Item("ProductLoanOutput")

$newTitle = "968-MS-EI-DB-DL-REC-NON-RNI-CTRFD-DL-MD-TR-1-LateRepayment"

# Update B1 text + style on both sheets.
# B22 on ProductLoanInput already carries the target style (fontId=1 Arial
# applied + green fill), so copy its formatting onto both B1 cells.
$wsIn.Range("B22").Copy()
$wsIn.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$wsIn.Range("B1").Value = $newTitle

$wsIn.Range("B22").Copy()
$wsOut.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$wsOut.Range("B1").Value = $newTitle

# Reset selections
$wsIn.Range("B1").Select()
$wsOut.Range("B1").Select()

# Switch active sheet to ProductLoanOutput
$wsOut.Activate()
